$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.646.07"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "2.020.15"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "247.79"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").Value = "0.636"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("D7").Value = "62.43"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +4.83%  "
$ws.Range("D10").Value = "57.79"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  +5.90%  "
$ws.Range("D12").Value = "0.103"
$ws.Range("E12").Value = "  -1.17%  "
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "23.23"
$ws.Range("E14").Value = "  +19.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.82%  "
$ws.Range("D16").Value = "2.315.83"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").Value = "5.53"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "2.023.14"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "36.580.35"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("D20").Value = "72.13"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "0.0₃0882"
$ws.Range("E21").Value = "  +2.89%  "
$ws.Range("D22").Value = "5.37"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").Value = "235.81"
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("E25").Value = "  -6.71%  "
$ws.Range("D26").Value = "2.32"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").Value = "9.81"
$ws.Range("E27").Value = "  +3.27%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "159.58"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "20.31"
$ws.Range("E29").Value = "  +3.62%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.133"
$ws.Range("E30").Value = "  +22.10%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "5.05"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").Value = "0.0619"
$ws.Range("E34").Value = "  +2.51%  "
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "2.37"
$ws.Range("E36").Value = "  -5.46%  "
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D37").Value = "6.36"
$ws.Range("E37").Value = "  +9.07%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").Value = "3.18"
$ws.Range("E40").Value = "  +27.26%  "
$ws.Range("D41").Value = "0.0996"
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("D42").Value = "1.24"
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("D45").Value = "16.98"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("D47").Value = "93.62"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "7.65"
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("D49").Value = "1.369.56"
$ws.Range("E49").Value = "  -3.62%  "
$ws.Range("D50").Value = "2.89"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "2.204.96"
$ws.Range("E51").Value = "  +0.22%  "
